# Actually evaluate open ended ranges.
#
# Sheet1 gains a header row: a new blank row is inserted above the
# existing data (old row 1, the literal 1/1/1 row, becomes row 2; the
# formula rows below it shift down to rows 3-5). The new row 1 is filled
# with text headers, and a new column D is added with a SUM() formula
# that evaluates the open-ended column A range.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 1 on Sheet1 - shifts the old data (and the
# cross-sheet formulas that point at it) down by one row, updating
# references automatically.
$ws1.Rows.Item(1).Insert() | Out-Null

# New header row.
$ws1.Range("A1").Value = "Hello"
$ws1.Range("B1").Value = "This "
$ws1.Range("C1").Value = "Is"
$ws1.Range("D1").Value = "A Table"

# New column D: evaluate the (now) open-ended column A range.
$ws1.Range("D2").Formula = "=SUM(Sheet1!A:A)"

# The workbook-level defined name pointed at the old last formula cell
# (C4); after the row insert that cell lives at C5.
$lastCell = $wb.Names.Item("LastCell")
$lastCell.RefersTo = "=Sheet1!`$C`$5"

# Active sheet/selection moved from Sheet2!D1 to Sheet1!M36.
$ws1.Activate() | Out-Null
$ws1.Range("M36").Select() | Out-Null
